# Apply crypto price/volume refresh (GitHub Actions scheduled data update).
# Generated from the cell-level diff between before.xlsx and after.xlsx.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = '35.388.98'
$ws.Range("E2").Value = '  -0.44%  '

# Row 3
$ws.Range("D3").Value = '1.917.53'
$ws.Range("E3").Value = '  +0.22%  '

# Row 4
$ws.Range("E4").Value = '  -0.35%  '

# Row 5
$ws.Range("D5").Value = '''0.719'
$ws.Range("E5").Value = '  +9.87%  '

# Row 6
$ws.Range("D6").Value = '''252.99'
$ws.Range("E6").Value = '  +2.80%  '

# Row 7
$ws.Range("E7").Value = '  -0.30%  '

# Row 8
$ws.Range("D8").Value = '''40.74'
$ws.Range("E8").Value = '  -3.07%  '

# Row 9
$ws.Range("D9").Value = '''0.357'
$ws.Range("E9").Value = '  +2.74%  '

# Row 10
$ws.Range("D10").Value = '''52.99'
$ws.Range("E10").Value = '  +7.15%  '

# Row 11
$ws.Range("D11").Value = '''0.0734'
$ws.Range("E11").Value = '  +2.01%  '

# Row 12
$ws.Range("D12").Value = '''0.0997'

# Row 13
$ws.Range("D13").Value = '2.196.58'
$ws.Range("E13").Value = '  +0.10%  '

# Row 14
$ws.Range("D14").Value = '''12.60'
$ws.Range("E14").Value = '  +2.52%  '

# Row 15
$ws.Range("D15").Value = '''0.717'
$ws.Range("E15").Value = '  +2.46%  '

# Row 16
$ws.Range("B16").Value = 'Polkadot'
$ws.Range("C16").Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range("D16").Value = '''4.92'
$ws.Range("E16").Value = '  +0.43%  '

# Row 17
$ws.Range("B17").Value = 'WrappedEther'
$ws.Range("C17").Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range("D17").Value = '1.917.14'
$ws.Range("E17").Value = '  +0.05%  '

# Row 18
$ws.Range("D18").Value = '35.427.25'
$ws.Range("E18").Value = '  -0.38%  '

# Row 19
$ws.Range("D19").Value = '''73.15'

# Row 20
$ws.Range("D20").Value = '0.0₃0828'
$ws.Range("E20").Value = '  +0.56%  '

# Row 21
$ws.Range("D21").Value = '''13.12'
$ws.Range("E21").Value = '  +3.68%  '

# Row 22
$ws.Range("D22").Value = '''241.83'
$ws.Range("E22").Value = '  -1.47%  '

# Row 23
$ws.Range("D23").Value = '''5.07'
$ws.Range("E23").Value = '  +4.81%  '

# Row 24
$ws.Range("E24").Value = '  -0.44%  '

# Row 25
$ws.Range("B25").Value = 'PancakeSwap'
$ws.Range("C25").Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range("D25").Value = '''2.39'
$ws.Range("E25").Value = '  +8.58%  '

# Row 26
$ws.Range("B26").Value = 'Toncoin'
$ws.Range("C26").Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range("D26").Value = '''2.32'
$ws.Range("E26").Value = '  +0.31%  '

# Row 27
$ws.Range("D27").Value = '''167.60'
$ws.Range("E27").Value = '  -2.15%  '

# Row 28
$ws.Range("D28").Value = '''8.70'
$ws.Range("E28").Value = '  +3.34%  '

# Row 29
$ws.Range("D29").Value = '''0.134'
$ws.Range("E29").Value = '  +5.00%  '

# Row 30
$ws.Range("D30").Value = '''18.73'
$ws.Range("E30").Value = '  +1.20%  '

# Row 31
$ws.Range("D31").Value = '4.132.52'

# Row 32
$ws.Range("E32").Value = '  +4.42%  '

# Row 33
$ws.Range("D33").Value = '''1.99'
$ws.Range("E33").Value = '  +13.27%  '

# Row 34
$ws.Range("D34").Value = '''0.0579'
$ws.Range("E34").Value = '  +1.49%  '

# Row 35
$ws.Range("D35").Value = '''4.27'
$ws.Range("E35").Value = '  +2.80%  '

# Row 36
$ws.Range("B36").Value = 'TrustWalletToken'
$ws.Range("C36").Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range("D36").Value = '''1.58'
$ws.Range("E36").Value = '  +16.95%  '

# Row 37
$ws.Range("B37").Value = 'BinanceUSD'
$ws.Range("C37").Value = 'https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd'
$ws.Range("D37").Value = '''1.01'
$ws.Range("E37").Value = '  -0.30%  '

# Row 38
$ws.Range("B38").Value = 'ImmutableX'
$ws.Range("C38").Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range("D38").Value = '''0.914'
$ws.Range("E38").Value = '  -0.79%  '

# Row 39
$ws.Range("D39").Value = '''2.08'
$ws.Range("E39").Value = '  +2.16%  '

# Row 40
$ws.Range("D40").Value = '''17.51'
$ws.Range("E40").Value = '  +11.05%  '

# Row 41
$ws.Range("D41").Value = '''98.12'
$ws.Range("E41").Value = '  +7.15%  '

# Row 42
$ws.Range("E42").Value = '  +2.90%  '

# Row 43
$ws.Range("D43").Value = '''0.0210'
$ws.Range("E43").Value = '  -0.99%  '

# Row 44
$ws.Range("D44").Value = '''0.0652'
$ws.Range("E44").Value = '  +1.50%  '

# Row 45
$ws.Range("D45").Value = '''2.50'
$ws.Range("E45").Value = '  +4.14%  '

# Row 46
$ws.Range("D46").Value = '1.346.25'
$ws.Range("E46").Value = '  -0.80%  '

# Row 47
$ws.Range("D47").Value = '''2.43'
$ws.Range("E47").Value = '  +1.05%  '

# Row 48
$ws.Range("D48").Value = '''2.78'
$ws.Range("E48").Value = '  -0.61%  '

# Row 49
$ws.Range("B49").Value = 'MultiversX'
$ws.Range("C49").Value = 'https://coinranking.com/coin/omwkOTglq+multiversx-egld'
$ws.Range("D49").Value = '''45.27'
$ws.Range("E49").Value = '  -5.39%  '

# Row 50
$ws.Range("B50").Value = 'FraxShare'
$ws.Range("C50").Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range("D50").Value = '''6.61'
$ws.Range("E50").Value = '  +0.51%  '

# Row 51
$ws.Range("D51").Value = '''11.98'
$ws.Range("E51").Value = '  -5.91%  '
